{"js": "// Update the answer key table: replace the arithmetic-problem text in each\n// populated data cell with its new value, while leaving cell/paragraph/run\n// formatting untouched. Cells are addressed by (row, col) in the single\n// table so that duplicate/overlapping text values (e.g. \"62\u00f74=15, 2\" is\n// both an old value in one cell and a new value in another) never cause\n// ambiguous matches.\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  return;\n}\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"41\u00f72=20, 1\", newText: \"33\u00f76=5, 3\" },\n  { row: 0, col: 1, oldText: \"14\u00f77=2, 0\", newText: \"94\u00f78=11, 6\" },\n  { row: 0, col: 2, oldText: \"27\u00f78=3, 3\", newText: \"40\u00f76=6, 4\" },\n  { row: 0, col: 3, oldText: \"33\u00f74=8, 1\", newText: \"88\u00f74=22, 0\" },\n  { row: 0, col: 4, oldText: \"48\u00f74=12, 0\", newText: \"90\u00f79=10, 0\" },\n  { row: 4, col: 0, oldText: \"51\u00f79=5, 6\", newText: \"89\u00f78=11, 1\" },\n  { row: 4, col: 1, oldText: \"88\u00f72=44, 0\", newText: \"54\u00f79=6, 0\" },\n  { row: 4, col: 2, oldText: \"30\u00f72=15, 0\", newText: \"10\u00f75=2, 0\" },\n  { row: 4, col: 3, oldText: \"36\u00f77=5, 1\", newText: \"90\u00f74=22, 2\" },\n  { row: 4, col: 4, oldText: \"27\u00f75=5, 2\", newText: \"69\u00f76=11, 3\" },\n  { row: 8, col: 0, oldText: \"33\u00f79=3, 6\", newText: \"14\u00f76=2, 2\" },\n  { row: 8, col: 1, oldText: \"30\u00f79=3, 3\", newText: \"68\u00f75=13, 3\" },\n  { row: 8, col: 2, oldText: \"58\u00f73=19, 1\", newText: \"62\u00f76=10, 2\" },\n  { row: 8, col: 3, oldText: \"75\u00f72=37, 1\", newText: \"79\u00f72=39, 1\" },\n  { row: 8, col: 4, oldText: \"51\u00f73=17, 0\", newText: \"17\u00f72=8, 1\" },\n  { row: 12, col: 0, oldText: \"25\u00f77=3, 4\", newText: \"90\u00f74=22, 2\" },\n  { row: 12, col: 1, oldText: \"12\u00f75=2, 2\", newText: \"72\u00f78=9, 0\" },\n  { row: 12, col: 2, oldText: \"62\u00f74=15, 2\", newText: \"58\u00f72=29, 0\" },\n  { row: 12, col: 3, oldText: \"22\u00f75=4, 2\", newText: \"33\u00f78=4, 1\" },\n  { row: 12, col: 4, oldText: \"39\u00f72=19, 1\", newText: \"45\u00f77=6, 3\" },\n  { row: 16, col: 0, oldText: \"30\u00f78=3, 6\", newText: \"81\u00f74=20, 1\" },\n  { row: 16, col: 2, oldText: \"77\u00f74=19, 1\", newText: \"62\u00f74=15, 2\" },\n  { row: 16, col: 3, oldText: \"20\u00f77=2, 6\", newText: \"97\u00f77=13, 6\" },\n  { row: 16, col: 4, oldText: \"81\u00f75=16, 1\", newText: \"20\u00f76=3, 2\" },\n];\n\n// Load each target cell's body paragraphs so we can grab a Range to\n// replace text in-place (keeps run/paragraph formatting intact).\nconst cells = replacements.map((r) => table.getCell(r.row, r.col));\ncells.forEach((cell) => cell.body.paragraphs.load(\"items/text\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const { oldText, newText } = replacements[i];\n  const paragraph = cells[i].body.paragraphs.items[0];\n  const currentText = paragraph.text;\n  // Only touch the paragraph if it still holds the expected original\n  // text; this keeps the script a no-op for any cell that doesn't match\n  // (defensive, in case the table shape ever differs from what we expect).\n  if (currentText === oldText) {\n    paragraph.getRange().insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the answer key table: replace the arithmetic-problem text in each\n# populated data cell with its new value, while leaving cell/paragraph/run\n# formatting untouched. Cells are addressed by (row, col) in the single\n# table so that duplicate/overlapping text values (e.g. \"62\u00f74=15, 2\" is\n# both an old value in one cell and a new value in another) never cause\n# ambiguous matches.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"41\u00f72=20, 1\"; NewText = \"33\u00f76=5, 3\" }\n    @{ Row = 1; Col = 2; OldText = \"14\u00f77=2, 0\"; NewText = \"94\u00f78=11, 6\" }\n    @{ Row = 1; Col = 3; OldText = \"27\u00f78=3, 3\"; NewText = \"40\u00f76=6, 4\" }\n    @{ Row = 1; Col = 4; OldText = \"33\u00f74=8, 1\"; NewText = \"88\u00f74=22, 0\" }\n    @{ Row = 1; Col = 5; OldText = \"48\u00f74=12, 0\"; NewText = \"90\u00f79=10, 0\" }\n    @{ Row = 5; Col = 1; OldText = \"51\u00f79=5, 6\"; NewText = \"89\u00f78=11, 1\" }\n    @{ Row = 5; Col = 2; OldText = \"88\u00f72=44, 0\"; NewText = \"54\u00f79=6, 0\" }\n    @{ Row = 5; Col = 3; OldText = \"30\u00f72=15, 0\"; NewText = \"10\u00f75=2, 0\" }\n    @{ Row = 5; Col = 4; OldText = \"36\u00f77=5, 1\"; NewText = \"90\u00f74=22, 2\" }\n    @{ Row = 5; Col = 5; OldText = \"27\u00f75=5, 2\"; NewText = \"69\u00f76=11, 3\" }\n    @{ Row = 9; Col = 1; OldText = \"33\u00f79=3, 6\"; NewText = \"14\u00f76=2, 2\" }\n    @{ Row = 9; Col = 2; OldText = \"30\u00f79=3, 3\"; NewText = \"68\u00f75=13, 3\" }\n    @{ Row = 9; Col = 3; OldText = \"58\u00f73=19, 1\"; NewText = \"62\u00f76=10, 2\" }\n    @{ Row = 9; Col = 4; OldText = \"75\u00f72=37, 1\"; NewText = \"79\u00f72=39, 1\" }\n    @{ Row = 9; Col = 5; OldText = \"51\u00f73=17, 0\"; NewText = \"17\u00f72=8, 1\" }\n    @{ Row = 13; Col = 1; OldText = \"25\u00f77=3, 4\"; NewText = \"90\u00f74=22, 2\" }\n    @{ Row = 13; Col = 2; OldText = \"12\u00f75=2, 2\"; NewText = \"72\u00f78=9, 0\" }\n    @{ Row = 13; Col = 3; OldText = \"62\u00f74=15, 2\"; NewText = \"58\u00f72=29, 0\" }\n    @{ Row = 13; Col = 4; OldText = \"22\u00f75=4, 2\"; NewText = \"33\u00f78=4, 1\" }\n    @{ Row = 13; Col = 5; OldText = \"39\u00f72=19, 1\"; NewText = \"45\u00f77=6, 3\" }\n    @{ Row = 17; Col = 1; OldText = \"30\u00f78=3, 6\"; NewText = \"81\u00f74=20, 1\" }\n    @{ Row = 17; Col = 3; OldText = \"77\u00f74=19, 1\"; NewText = \"62\u00f74=15, 2\" }\n    @{ Row = 17; Col = 4; OldText = \"20\u00f77=2, 6\"; NewText = \"97\u00f77=13, 6\" }\n    @{ Row = 17; Col = 5; OldText = \"81\u00f75=16, 1\"; NewText = \"20\u00f76=3, 2\" }\n)\n\nforeach ($item in $replacements) {\n    $cell = $table.Cell($item.Row, $item.Col)\n    $currentText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    # Only touch the cell if it still holds the expected original text;\n    # this keeps the script a no-op for any cell that doesn't match\n    # (defensive, in case the table shape ever differs from what we expect).\n    if ($currentText -eq $item.OldText) {\n        $cell.Range.Text = $item.NewText\n    }\n}\n"}
